$d = $word.ActiveDocument

# The cover letter referred to the "Drug Design Data Resource (workshop)";
# correct the parenthetical to the resource's actual acronym, (D3R).
$d.Content.Find.Execute(
    "(workshop)",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "(D3R)",
    2
)
